$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C for "Ma danh muc" (category code),
# shifting the old C:F (Ton dau ky / Nhap / Xuat / Ton cuoi ky) into D:G.
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = 9.2

# Restore the period label in C1 (it was pushed to D1 by the insert) and
# re-merge the title bar across the now-wider C1:G1 range.
$ws.Range("C1").Value = "Thời gian:  06_2024"
$ws.Range("C1:G1").Merge()

# New column header
$ws.Range("C2").Value = "Mã danh mục"

# Fill in the new "Mã danh mục" values for the existing product rows
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1

# Append the new product row
$ws.Range("A9").Value = 34
$ws.Range("B9").Value = "may giat 1"
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 5
